# -refactoring -using webeye framework to detect available cameras
#
# Diary sheet update: correct the end time logged for 16.11.18 (row 26)
# and log a new entry for 19.11.18 (row 27), then move the active
# selection to B28 to reflect where data entry continues.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 (16.11.18): end time corrected -> duration recalculates via the
# existing shared formula in column D.
$ws.Range("C26").Value = 0.59375

# Row 27: new diary entry for 19.11.18, same "-refactoring" remark as the
# previous day.
$ws.Range("A27").Value = "19.11.18"
$ws.Range("B27").Value = 0.5
$ws.Range("C27").Value = 0.79166666666666663
$ws.Range("D27").Formula = "=C27-B27"
$ws.Range("E27").Value = "-refactoring"

# Recalculate so the total in D36 and the new D27 duration are up to date.
$excel.CalculateFull()

# Reflect where the user's selection ended up after the edit.
$ws.Range("B28").Select()
